# FIX — Lidando com valores None nas métricas. Computando métricas para
# resultados sem otimização. Updates the descriptive-statistics values in
# the active worksheet (compas_stats_carla.xlsx) to reflect the
# recomputed metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - mean
$ws.Range("B3").Value = 5.161616161616162
$ws.Range("C3").Value = 2.199146899475082
$ws.Range("D3").Value = 2.163708741019674
$ws.Range("E3").Value = 1

# Row 4 - std
$ws.Range("B4").Value = 0.3699716124540519
$ws.Range("C4").Value = 0.4004557400923186
$ws.Range("D4").Value = 0.3737509168577619
$ws.Range("E4").Value = 0

# Row 5 - min
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 2.003223017231418
$ws.Range("D5").Value = 2.000004837570589
$ws.Range("E5").Value = 1

# Row 6 - 25%
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 2.015855129061242
$ws.Range("D6").Value = 2.000111827743817
$ws.Range("E6").Value = 1

# Row 7 - 50%
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 2.022797494580614
$ws.Range("D7").Value = 2.000247305107726
$ws.Range("E7").Value = 1

# Row 8 - 75%
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 2.04142350414401
$ws.Range("D8").Value = 2.000910619094734
$ws.Range("E8").Value = 1

# Row 9 - max
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 3.306349560875093
$ws.Range("D9").Value = 3.041236137226413
$ws.Range("E9").Value = 1
